# The document contains the bold phrase "DOCX, DOC, PDF, HTML, XPS, R" followed
# by a "_GoBack" bookmark and then the bold phrase "TF and TXT" (two separate
# runs split around the bookmark). Word's Find/Replace treats bookmarks as
# transparent, so searching for the full joined phrase "DOCX, DOC, PDF, HTML,
# XPS, RTF and TXT" finds it across both runs; replacing it in-place merges
# the text into a single run (dropping the now-redundant bookmark) while
# keeping the surrounding bold formatting.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",  # FindText
    $false,                                    # MatchCase
    $false,                                    # MatchWholeWord
    $false,                                    # MatchWildcards
    $false,                                    # MatchSoundsLike
    $false,                                    # MatchAllWordForms
    $true,                                     # Forward
    1,                                         # Wrap (wdFindContinue)
    $false,                                    # Format
    "DOCX, DOC, PDF, HTML, XPS, RTF and TXT",  # ReplaceWith
    2                                          # Replace (wdReplaceAll)
)
